# Manual Testing Log - replace "Infinity" divide-by-zero tests with
# proper "Cannot Divide By Zero" error-handling tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Log")

# Remove the obsolete "[Divide by Zero] result of Infinity from: 7 / 0" row.
$ws.Rows.Item(16).Delete()

# Remove the obsolete "[Negative] equation does NOT contain "Infinity" ..." row
# (this row is now at index 21 after the previous delete shifted rows up).
$ws.Rows.Item(21).Delete()

# Append the two new Divide-By-Zero regression tests at the bottom of the log.
$ws.Cells.Item(45, 1).Value = "2020-30-09"
$ws.Cells.Item(45, 2).Value = '[Divide By Zero] "Cannot Divide By Zero" result from: 9 / 0'
$ws.Cells.Item(45, 3).Value = "Y"

$ws.Cells.Item(46, 1).Value = "2020-30-09"
$ws.Cells.Item(46, 2).Value = '[Divide By Zero] "Cannot Divide By Zero" result from: 12 * 8 / 0'
$ws.Cells.Item(46, 3).Value = "Y"

# The old backspace-on-negative test used to report a lone "-" when the
# result hit zero; it was fixed to show "0" instead.
$ws.Cells.Item(41, 2).Value = '[Backspace] result shows "0" from: -55, backspace, backspace'

# Move the active selection like the author left it.
$ws.Range("F25").Select()
